# Hortaliza, Agrícola del Norte S.A. de Arica - Cebolla
# Commit: "Fruta / hortaliza, semanal"
#
# A new weekly price entry (3 rows: 1a/2a/3a cosecha, variedad "Sin especificar",
# fecha 2023-04-05) is inserted at the top of the Agrícola del Norte / Cebolla
# block (row 1201), pushing the existing 1201:1296 rows down to 1204:1299.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at 1201, shifting everything from 1201 downward by 3.
$ws.Rows("1201:1203").Insert()

# --- Row 1201 : 1a (cosecha) ---
$ws.Cells.Item(1201,1).Value  = 1
$ws.Cells.Item(1201,2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(1201,3).Value  = "Arica y Parinacota"
$ws.Cells.Item(1201,4).Value  = 45021
$ws.Cells.Item(1201,5).Value  = 15
$ws.Cells.Item(1201,6).Value  = 100112004
$ws.Cells.Item(1201,7).Value  = "Cebolla"
$ws.Cells.Item(1201,8).Value  = "Sin especificar"
$ws.Cells.Item(1201,9).Value  = "1a (cosecha)"
$ws.Cells.Item(1201,10).Value = 250
$ws.Cells.Item(1201,11).Value = 9000
$ws.Cells.Item(1201,12).Value = 10000
$ws.Cells.Item(1201,13).Value = 9500
$ws.Cells.Item(1201,14).Value = "`$/malla 18 kilos"
$ws.Cells.Item(1201,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(1201,16).Value = 528
$ws.Cells.Item(1201,17).Value = 18
$ws.Cells.Item(1201,18).Value = "Hortaliza"

# --- Row 1202 : 2a (cosecha) ---
$ws.Cells.Item(1202,1).Value  = 1
$ws.Cells.Item(1202,2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(1202,3).Value  = "Arica y Parinacota"
$ws.Cells.Item(1202,4).Value  = 45021
$ws.Cells.Item(1202,5).Value  = 15
$ws.Cells.Item(1202,6).Value  = 100112004
$ws.Cells.Item(1202,7).Value  = "Cebolla"
$ws.Cells.Item(1202,8).Value  = "Sin especificar"
$ws.Cells.Item(1202,9).Value  = "2a (cosecha)"
$ws.Cells.Item(1202,10).Value = 300
$ws.Cells.Item(1202,11).Value = 7000
$ws.Cells.Item(1202,12).Value = 8000
$ws.Cells.Item(1202,13).Value = 7500
$ws.Cells.Item(1202,14).Value = "`$/malla 18 kilos"
$ws.Cells.Item(1202,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(1202,16).Value = 417
$ws.Cells.Item(1202,17).Value = 18
$ws.Cells.Item(1202,18).Value = "Hortaliza"

# --- Row 1203 : 3a (cosecha) ---
$ws.Cells.Item(1203,1).Value  = 1
$ws.Cells.Item(1203,2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(1203,3).Value  = "Arica y Parinacota"
$ws.Cells.Item(1203,4).Value  = 45021
$ws.Cells.Item(1203,5).Value  = 15
$ws.Cells.Item(1203,6).Value  = 100112004
$ws.Cells.Item(1203,7).Value  = "Cebolla"
$ws.Cells.Item(1203,8).Value  = "Sin especificar"
$ws.Cells.Item(1203,9).Value  = "3a (cosecha)"
$ws.Cells.Item(1203,10).Value = 300
$ws.Cells.Item(1203,11).Value = 5000
$ws.Cells.Item(1203,12).Value = 6000
$ws.Cells.Item(1203,13).Value = 5500
$ws.Cells.Item(1203,14).Value = "`$/malla 18 kilos"
$ws.Cells.Item(1203,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(1203,16).Value = 306
$ws.Cells.Item(1203,17).Value = 18
$ws.Cells.Item(1203,18).Value = "Hortaliza"
